# Apply the "updated front-end data" edit to the All_Jockeys sheet:
# replace the jockey "Dylan Davis" (row 5) with new jockey "Andre Worrie",
# keeping only the columns that have data for the new jockey (weights),
# and clearing out the columns for which no data is (yet) available.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: swap Dylan Davis -> Andre Worrie with new min/max/avg weight values.
$ws.Range("A5").Value = "Andre Worrie"
$ws.Range("B5").Value = 118
$ws.Range("C5").Value = 126
$ws.Range("D5").Value = 121.686055726376

# Height / Age / years_of_experience are unknown for the new jockey.
$ws.Range("E5:G5").ClearContents()

# 2023/all-time firsts/seconds/thirds are unknown for the new jockey
# (cells remain present, keeping their number formatting, but empty).
$ws.Range("H5:M5").ClearContents()

# No profile image link for the new jockey.
$ws.Range("N5").ClearContents()

# Update the last-selected cell, as recorded by Excel when the file was saved.
$ws.Range("F15").Select()
